$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:D3").Value = -0.04
$ws.Range("E2:E3").Value = 0.0669
$ws.Range("F2:F3").Value = -0.00722
$ws.Range("G2:G3").Value = 0.1397621963148399
$ws.Range("H2:H3").Value = 0.1397621963148399
$ws.Range("I2:I3").Value = 0.1388130103685885
$ws.Range("J2:J3").Value = 0.1142119803225112
$ws.Range("K2:K3").Value = 1286.5
$ws.Range("L2:L3").Value = 0.09321585647728838
$ws.Range("U2:U3").Value = 2604.2
$ws.Range("V2:V3").Value = 0.2612193311533292
$ws.Range("W2:W3").Value = 0.1051242451727829
$ws.Range("X2:X3").Value = 0.06126085697584728
$ws.Range("Y2:Y3").Value = 0.04386338819693564
$ws.Range("Z2:Z3").Value = 0.9446798316164139
$ws.Range("AA2:AA3").Value = 0.1078937543396471
$ws.Range("AB2:AB3").Value = 0.04376391170335094
$ws.Range("AC2:AC3").Value = 0.06412984263629615
$ws.Range("AD2:AD3").Value = 6765.5
$ws.Range("AE2:AE3").Value = 0
$ws.Range("AF2:AF3").Value = 6765.5
$ws.Range("AG2:AG3").Value = 4161.3
$ws.Range("AH2:AH3").Value = 0.4042748985652737
$ws.Range("AI2:AI3").Value = 0.2987226301544059
$ws.Range("AJ2:AJ3").Value = 0.2944864727154352
$ws.Range("AK2:AK3").Value = 0.207609297591786
$ws.Range("AL2:AL3").Value = 117.5
$ws.Range("AM2:AM3").Value = 117.5
$ws.Range("AN2:AN3").Value = 2.950115554005145
$ws.Range("AO2:AO3").Value = 16.30468085106383
$ws.Range("AP2:AP3").Value = 1.814546723062835
$ws.Range("AQ2:AQ3").Value = 16.30468085106383
